$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 6944
$ws.Cells.Item(3, 9).Value = 7197
$ws.Cells.Item(4, 9).Value = 1652
$ws.Cells.Item(5, 9).Value = 676
$ws.Cells.Item(6, 9).Value = 8500
$ws.Cells.Item(7, 9).Value = 24969

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 9).Value = 120
$ws.Cells.Item(7, 9).Value = 292

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(4, 9).Value = 6
$ws.Cells.Item(7, 9).Value = 141

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 9).Value = 250
$ws.Cells.Item(6, 9).Value = 235
$ws.Cells.Item(7, 9).Value = 771

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(6, 9).Value = 113
$ws.Cells.Item(7, 9).Value = 436

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 225
$ws.Cells.Item(3, 9).Value = 353
$ws.Cells.Item(7, 9).Value = 945

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 9).Value = 174
$ws.Cells.Item(7, 9).Value = 581

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 9).Value = 198
$ws.Cells.Item(5, 9).Value = 75
$ws.Cells.Item(6, 9).Value = 181
$ws.Cells.Item(7, 9).Value = 781
$ws.Cells.Item(9, 9).Value = 131
$ws.Cells.Item(10, 9).Value = 182
$ws.Cells.Item(14, 9).Value = 141
$ws.Cells.Item(15, 9).Value = 291
$ws.Cells.Item(16, 9).Value = 73
$ws.Cells.Item(17, 9).Value = 37
$ws.Cells.Item(18, 9).Value = 195
$ws.Cells.Item(19, 9).Value = 700
$ws.Cells.Item(20, 9).Value = 618
$ws.Cells.Item(23, 9).Value = 245
$ws.Cells.Item(26, 9).Value = 36
$ws.Cells.Item(33, 9).Value = 1098
$ws.Cells.Item(36, 9).Value = 339
$ws.Cells.Item(37, 9).Value = 771
$ws.Cells.Item(42, 9).Value = 939
$ws.Cells.Item(43, 9).Value = 213
$ws.Cells.Item(44, 9).Value = 189
$ws.Cells.Item(48, 9).Value = 318
$ws.Cells.Item(50, 9).Value = 131
$ws.Cells.Item(52, 9).Value = 565
$ws.Cells.Item(54, 9).Value = 490
$ws.Cells.Item(55, 9).Value = 288
$ws.Cells.Item(59, 9).Value = 42
$ws.Cells.Item(63, 9).Value = 75
$ws.Cells.Item(65, 9).Value = 581
$ws.Cells.Item(66, 9).Value = 74
$ws.Cells.Item(67, 9).Value = 945
$ws.Cells.Item(73, 9).Value = 227
$ws.Cells.Item(78, 9).Value = 333
$ws.Cells.Item(79, 9).Value = 715
$ws.Cells.Item(83, 9).Value = 535
$ws.Cells.Item(85, 9).Value = 1112
$ws.Cells.Item(86, 9).Value = 161
$ws.Cells.Item(87, 9).Value = 65
$ws.Cells.Item(88, 9).Value = 230
$ws.Cells.Item(90, 9).Value = 325
$ws.Cells.Item(94, 9).Value = 254
$ws.Cells.Item(95, 9).Value = 384
$ws.Cells.Item(96, 9).Value = 292
$ws.Cells.Item(97, 9).Value = 226
$ws.Cells.Item(99, 9).Value = 436
$ws.Cells.Item(101, 9).Value = 24969

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 178
$ws.Cells.Item(7, 9).Value = 535

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(2, 9).Value = 137
$ws.Cells.Item(3, 9).Value = 133
$ws.Cells.Item(7, 9).Value = 384

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(5, 9).Value = 46
$ws.Cells.Item(6, 9).Value = 353
$ws.Cells.Item(7, 9).Value = 1098

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 9).Value = 106
$ws.Cells.Item(6, 9).Value = 236
$ws.Cells.Item(7, 9).Value = 490

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 9).Value = 204
$ws.Cells.Item(7, 9).Value = 700

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 9).Value = 63
$ws.Cells.Item(6, 9).Value = 54
$ws.Cells.Item(7, 9).Value = 189

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 9).Value = 52
$ws.Cells.Item(3, 9).Value = 60
$ws.Cells.Item(6, 9).Value = 163
$ws.Cells.Item(7, 9).Value = 318

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 9).Value = 418
$ws.Cells.Item(7, 9).Value = 1112

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(6, 9).Value = 53
$ws.Cells.Item(7, 9).Value = 181

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 9).Value = 377
$ws.Cells.Item(7, 9).Value = 939

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(6, 9).Value = 85
$ws.Cells.Item(7, 9).Value = 182

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(3, 9).Value = 86
$ws.Cells.Item(7, 9).Value = 333

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 9).Value = 87
$ws.Cells.Item(7, 9).Value = 288

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 9).Value = 67
$ws.Cells.Item(7, 9).Value = 245

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(3, 9).Value = 235
$ws.Cells.Item(7, 9).Value = 715

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(6, 9).Value = 217
$ws.Cells.Item(7, 9).Value = 618

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 9).Value = 90
$ws.Cells.Item(7, 9).Value = 195

$ws = $wb.Worksheets.Item('Burnside')
$ws.Cells.Item(2, 9).Value = 11
$ws.Cells.Item(7, 9).Value = 37

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(6, 9).Value = 106
$ws.Cells.Item(7, 9).Value = 339

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 140
$ws.Cells.Item(4, 9).Value = 39
$ws.Cells.Item(6, 9).Value = 184
$ws.Cells.Item(7, 9).Value = 565

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 9).Value = 50
$ws.Cells.Item(7, 9).Value = 254

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(6, 9).Value = 111
$ws.Cells.Item(7, 9).Value = 291

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 9).Value = 41
$ws.Cells.Item(7, 9).Value = 131

$ws = $wb.Worksheets.Item('East Village')
$ws.Cells.Item(6, 9).Value = 22
$ws.Cells.Item(7, 9).Value = 36

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 9).Value = 33
$ws.Cells.Item(7, 9).Value = 74

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 9).Value = 46
$ws.Cells.Item(7, 9).Value = 131

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(6, 9).Value = 61
$ws.Cells.Item(7, 9).Value = 227

$ws = $wb.Worksheets.Item('Montclare')
$ws.Cells.Item(6, 9).Value = 14
$ws.Cells.Item(7, 9).Value = 42

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 9).Value = 42
$ws.Cells.Item(7, 9).Value = 198

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(2, 9).Value = 37
$ws.Cells.Item(6, 9).Value = 147
$ws.Cells.Item(7, 9).Value = 226

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(6, 9).Value = 73
$ws.Cells.Item(7, 9).Value = 230

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 9).Value = 427
$ws.Cells.Item(6, 9).Value = 480

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(6, 9).Value = 30
$ws.Cells.Item(7, 9).Value = 75

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 9).Value = 78
$ws.Cells.Item(7, 9).Value = 161

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 9).Value = 116
$ws.Cells.Item(7, 9).Value = 325

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(2, 9).Value = 43
$ws.Cells.Item(7, 9).Value = 213

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 9).Value = 256
$ws.Cells.Item(6, 9).Value = 212
$ws.Cells.Item(7, 9).Value = 781

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 9).Value = 39
$ws.Cells.Item(7, 9).Value = 65

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Cells.Item(6, 9).Value = 50
$ws.Cells.Item(7, 9).Value = 73
